$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel, so they remain text like the source data.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = '46.353.77'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").Value = '2.615.02'
$ws.Range("E3").Value = '  +3.88%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '307.96'
$ws.Range("E5").Value = '  +2.83%  '
$ws.Range("D6").Value = '100.37'
$ws.Range("E6").Value = '  +2.39%  '
$ws.Range("E7").Value = '  +2.62%  '
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("D9").Value = '0.578'
$ws.Range("E9").Value = '  +6.89%  '
$ws.Range("D10").Value = '39.55'
$ws.Range("E10").Value = '  +8.18%  '
$ws.Range("D11").Value = '0.0846'
$ws.Range("E11").Value = '  +5.93%  '
$ws.Range("D12").Value = '54.15'
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").Value = '8.17'
$ws.Range("E13").Value = '  +8.47%  '
$ws.Range("D14").Value = '3.006.83'
$ws.Range("E14").Value = '  +3.47%  '
$ws.Range("E15").Value = '  +1.33%  '
$ws.Range("D16").Value = '2.611.89'
$ws.Range("E16").Value = '  +4.25%  '
$ws.Range("D17").Value = '0.921'
$ws.Range("E17").Value = '  +4.37%  '
$ws.Range("D18").Value = '15.00'
$ws.Range("E18").Value = '  +2.94%  '
$ws.Range("D19").Value = '46.483.65'
$ws.Range("E19").Value = '  +0.36%  '
$ws.Range("E20").Value = '  +4.70%  '
$ws.Range("D21").Value = '12.96'
$ws.Range("E21").Value = '  -2.68%  '
$ws.Range("D22").Value = '6.75'
$ws.Range("E22").Value = '  +4.33%  '
$ws.Range("D23").Value = '71.61'
$ws.Range("E23").Value = '  +3.87%  '
$ws.Range("D24").Value = '274.62'
$ws.Range("E24").Value = '  +9.70%  '
$ws.Range("D25").Value = '3.04'
$ws.Range("E25").Value = '  +5.97%  '
$ws.Range("D26").Value = '2.17'
$ws.Range("E26").Value = '  +6.62%  '
$ws.Range("D27").Value = '28.88'
$ws.Range("E27").Value = '  +28.14%  '
$ws.Range("E28").Value = '  -0.26%  '
$ws.Range("E29").Value = '  -0.54%  '
$ws.Range("D30").Value = '10.64'
$ws.Range("E30").Value = '  +5.60%  '
$ws.Range("D31").Value = '2.29'
$ws.Range("E31").Value = '  +2.32%  '
$ws.Range("D32").Value = '39.06'
$ws.Range("E32").Value = '  -5.50%  '
$ws.Range("D33").Value = '6.39'
$ws.Range("E33").Value = '  +10.63%  '
$ws.Range("D34").Value = '3.66'
$ws.Range("E34").Value = '  -7.30%  '
$ws.Range("E35").Value = '  +2.24%  '
$ws.Range("D36").Value = '2.24'
$ws.Range("E36").Value = '  +1.35%  '
$ws.Range("D37").Value = '0.0842'
$ws.Range("E37").Value = '  +5.98%  '
$ws.Range("D38").Value = '151.02'
$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("E39").Value = '  +3.98%  '
$ws.Range("E40").Value = '  +4.92%  '
$ws.Range("D41").Value = '23.36'
$ws.Range("E41").Value = '  +39.26%  '
$ws.Range("D42").Value = '15.88'
$ws.Range("E42").Value = '  -0.20%  '
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").Value = '3.65'
$ws.Range("E43").Value = '  +8.04%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0332'
$ws.Range("E44").Value = '  +6.89%  '
$ws.Range("E45").Value = '  -1.49%  '
$ws.Range("D46").Value = '2.134.48'
$ws.Range("E46").Value = '  +4.81%  '
$ws.Range("E47").Value = '  -0.08%  '
$ws.Range("D48").Value = '93.34'
$ws.Range("E48").Value = '  +1.95%  '
$ws.Range("D49").Value = '9.51'
$ws.Range("E49").Value = '  +8.50%  '
$ws.Range("D50").Value = '1.79'
$ws.Range("E50").Value = '  -2.24%  '
$ws.Range("D51").Value = '109.33'
$ws.Range("E51").Value = '  +3.04%  '
